$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @("闻泰科技","方正科技","平潭发展"),
    @("平潭发展","东方财富","超颖电子"),
    @("香农芯创","闻泰科技","闻泰科技"),
    @("方正科技","香农芯创","方正科技"),
    @("中国卫星","平潭发展","盈新发展"),
    @("达华智能","中国卫星","大有能源"),
    @("格尔软件","国轩高科","山子高科"),
    @("寒武纪-U","中国核建","大洋电机"),
    @("国轩高科","贵州茅台","达华智能"),
    @("中国核建","达华智能","三花智控"),
    @("中际旭创","中信证券","中际旭创"),
    @("和而泰","中际旭创","和而泰"),
    @("江波龙","和而泰","神开股份"),
    @("东方财富","大为股份","上海电力"),
    @("中信证券","楚江新材","华建集团"),
    @("航天科技","江波龙","中国核建"),
    @("大为股份","格尔软件","国轩高科"),
    @("大洋电机","航天科技","中信证券"),
    @("三花智控","上海电力","香农芯创"),
    @("上海电力","寒武纪-U","楚江新材")
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowValues = $values[$i]
    $ws.Cells.Item($row, 1).Value = $rowValues[0]
    $ws.Cells.Item($row, 2).Value = $rowValues[1]
    $ws.Cells.Item($row, 3).Value = $rowValues[2]
}
